$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace header labels with numeric column indices (0-based)
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Row 2: shift the old row-1 header labels down into row 2
$ws.Range("A2").Value = "Lg.,mm"
$ws.Range("B2").Value = "Threading"
$ws.Range("C2").Value = "HeadDia., mm"
$ws.Range("D2").Value = "HeadHt., mm"
$ws.Range("E2").Value = "DriveSize"
$ws.Range("F2").Value = "TensileStrength, psi"
$ws.Range("G2").Value = "SpecificationsMet"
$ws.Range("H2").Value = "Pkg.Qty."
$ws.Range("J2").Value = "Pkg."

# Rows 3-26: set column L to "Grade 2 Titanium"
for ($row = 3; $row -le 26; $row++) {
    $ws.Cells.Item($row, 12).Value = "Grade 2 Titanium"
}
